# correcao do modo de calcular o estoque de produtos
$wb = $excel.ActiveWorkbook

# --- Sheet "produtos": fill in the product rows (A2:E16) ---
$wsProdutos = $wb.Worksheets.Item("produtos")

$produtos = @(
    @(1, "prof",          "0",  "KG",    3),
    @(2, "rd",             "0",  "kg",    0),
    @(3, "rr",             "33", "LITRO", 3),
    @(4, "ddf",            "2",  "KG",    0),
    @(5, "dew",            "2",  "LITRO", 3),
    @(6, "faw",            "2",  "LITRO", 0),
    @(7, "gseg",           "2",  "kg",    0),
    @(8, "esefse",         "2",  "LITRO", 3),
    @(9, "fse\ef",         "2",  "LITRO", 0),
    @(10, "fs\ef",         "2",  "KG",    0),
    @(11, "\ef\s",         "2",  "LITRO", 0),
    @(12, "sefs",          "2",  "LITRO", 0),
    @(13, "sf\es",         "2",  "LITRO", 0),
    @(14, "\efse",         "2",  "LITRO", 0),
    @(15, "\fefsdf\ese\e", "2",  "LITRO", 0)
)

# Column C stores the quantity as plain text (matches the source data),
# so format it as Text before writing the values.
$wsProdutos.Range("C2:C16").NumberFormat = "@"

$r = 2
foreach ($row in $produtos) {
    $wsProdutos.Cells.Item($r, 1).Value = $row[0]
    $wsProdutos.Cells.Item($r, 2).Value = $row[1]
    $wsProdutos.Cells.Item($r, 3).Value = $row[2]
    $wsProdutos.Cells.Item($r, 4).Value = $row[3]
    $wsProdutos.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# --- Sheet "movimentos": fill in the movement rows (A2:F5) ---
$wsMovimentos = $wb.Worksheets.Item("movimentos")

$movimentos = @(
    @(1, 8, "esefse", "ENTRADA", 3, "2026-02-26 16:46:08"),
    @(2, 3, "rr",     "ENTRADA", 3, "2026-02-26 16:48:25"),
    @(3, 5, "dew",    "ENTRADA", 3, "2026-02-26 17:01:20"),
    @(4, 1, "prof",   "ENTRADA", 3, "2026-02-26 17:01:21")
)

$r = 2
foreach ($row in $movimentos) {
    $wsMovimentos.Cells.Item($r, 1).Value = $row[0]
    $wsMovimentos.Cells.Item($r, 2).Value = $row[1]
    $wsMovimentos.Cells.Item($r, 3).Value = $row[2]
    $wsMovimentos.Cells.Item($r, 4).Value = $row[3]
    $wsMovimentos.Cells.Item($r, 5).Value = $row[4]
    $wsMovimentos.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
